# Add the first two pages of report template
# Rename the placeholder test values in column B to the "summer" test names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Me_summer_test1"
$ws.Range("B3").Value = "Me_summer_test2"

# Update the view: clear the frozen/scrolled top-left cell and move the
# active selection to B3 (matching the saved view state in the workbook).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B3").Select()
